$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.259.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "'1.676.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "'212.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.95%  "
$ws.Range("D6").Value = "'0.5280"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.62%  "
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").Value = "'0.2660"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").Value = "'0.06296"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").Value = "'0.07560"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "'1.722.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "'4.474"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").Value = "'0.5645"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").Value = "'67.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").Value = "'0.000008051"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("D17").Value = "'26.061.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "'4.838"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'188.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "'10.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.99%  "
$ws.Range("D22").Value = "'6.225"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "'150.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "'0.1257"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.83%  "
$ws.Range("D26").Value = "'7.611"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("D27").Value = "'16.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "'0.06242"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "'1.364"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("D31").Value = "'3.523"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("D32").Value = "'3.443"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.40%  "
$ws.Range("D33").Value = "'1.642"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.25%  "
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("D35").Value = "'0.6081"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "'2.738"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01624"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'6.116"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'1.103.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "'0.8723"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("D43").Value = "'100.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Value = "'1.828.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "'56.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'8.023"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("D49").Value = "'0.05234"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").Value = "'0.4256"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "'5.995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.42%  "
